$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = -0.01826323110202846
$ws.Cells.Item(2, 3).Value = 0.9826393846830302
$ws.Cells.Item(2, 4).Value = -4.39275848435103
$ws.Cells.Item(2, 5).Value = 2.042986564144792
$ws.Cells.Item(2, 6).Value = -0.7114458636743382
$ws.Cells.Item(2, 7).Value = 2.225589065839782
$ws.Cells.Item(2, 8).Value = 0.0007645290829910302
$ws.Cells.Item(2, 9).Value = 0.436410659880665
$ws.Cells.Item(2, 10).Value = 1.130630940556499
$ws.Cells.Item(2, 11).Value = [double]"1.60277389766494e-07"
$ws.Cells.Item(2, 12).Value = 0.002803927456123335
$ws.Cells.Item(2, 13).Value = 24.0317582440992
$ws.Cells.Item(2, 14).Value = 23.42966264507794
$ws.Cells.Item(2, 15).Value = 24.64932649912292
$ws.Cells.Item(3, 2).Value = -0.03029595207641897
$ws.Cells.Item(3, 3).Value = 1.383004879170442
$ws.Cells.Item(3, 4).Value = -5.640408279622426
$ws.Cells.Item(3, 5).Value = 7.125165549156054
$ws.Cells.Item(3, 6).Value = 0.4679121322862519
$ws.Cells.Item(3, 7).Value = 7.225593029269941
$ws.Cells.Item(3, 8).Value = 0.05322541909867877
$ws.Cells.Item(3, 9).Value = [double]"5.03724395528593e-10"
$ws.Cells.Item(3, 10).Value = 0.9520357280252342
$ws.Cells.Item(3, 11).Value = 0.01442972973850082
$ws.Cells.Item(3, 12).Value = [double]"5.667026374556035e-07"
$ws.Cells.Item(3, 13).Value = 11.43144901117667
$ws.Cells.Item(3, 14).Value = 10.52108915524121
$ws.Cells.Item(3, 15).Value = 12.42057971061229
$ws.Cells.Item(4, 2).Value = -0.04017864546916282
$ws.Cells.Item(4, 3).Value = 1.532656360878771
$ws.Cells.Item(4, 4).Value = -3.876760707166097
$ws.Cells.Item(4, 5).Value = 5.835333761873049
$ws.Cells.Item(4, 6).Value = 0.3873319885409055
$ws.Cells.Item(4, 7).Value = 1.32202219557075
$ws.Cells.Item(4, 8).Value = 0.0007875194029758885
$ws.Cells.Item(4, 9).Value = [double]"4.440892098500626e-16"
$ws.Cells.Item(4, 10).Value = 0.8626148382101905
$ws.Cells.Item(4, 11).Value = [double]"8.605465636455148e-10"
$ws.Cells.Item(4, 12).Value = 0.03742614294781015
$ws.Cells.Item(4, 13).Value = 5.365808670784342
$ws.Cells.Item(4, 14).Value = 4.980596159310264
$ws.Cells.Item(4, 15).Value = 5.78081453916787
$ws.Cells.Item(5, 2).Value = -0.02325140714297254
$ws.Cells.Item(5, 3).Value = 1.121189232478837
$ws.Cells.Item(5, 4).Value = -2.640608132136125
$ws.Cells.Item(5, 5).Value = 2.939929088187708
$ws.Cells.Item(5, 6).Value = 0.1183850693177807
$ws.Cells.Item(5, 7).Value = -0.5626259227735115
$ws.Cells.Item(5, 8).Value = 0.03903027297625253
$ws.Cells.Item(5, 9).Value = 0.009186125291129499
$ws.Cells.Item(5, 10).Value = 0.9880803789329855
$ws.Cells.Item(5, 11).Value = [double]"6.380927099471565e-10"
$ws.Cells.Item(5, 12).Value = 0.2121330002904008
$ws.Cells.Item(5, 13).Value = 5.585501576877372
$ws.Cells.Item(5, 14).Value = 5.35176341126747
$ws.Cells.Item(5, 15).Value = 5.829448252442639
$ws.Cells.Item(6, 2).Value = -0.03544880984631753
$ws.Cells.Item(6, 3).Value = 1.254641566490597
$ws.Cells.Item(6, 4).Value = -4.600628468372605
$ws.Cells.Item(6, 5).Value = 2.794978009296771
$ws.Cells.Item(6, 6).Value = -0.2267668131688624
$ws.Cells.Item(6, 7).Value = 0.8056087046200622
$ws.Cells.Item(6, 8).Value = 0.3039261120476587
$ws.Cells.Item(6, 9).Value = [double]"2.489870643807013e-05"
$ws.Cells.Item(6, 10).Value = 0.8181963584906322
$ws.Cells.Item(6, 11).Value = [double]"1.026204605077818e-22"
$ws.Cells.Item(6, 12).Value = 0.3028605377448993
$ws.Cells.Item(6, 13).Value = 10.66177127902857
$ws.Cells.Item(6, 14).Value = 9.840680809547653
$ws.Cells.Item(6, 15).Value = 11.55137220750316
$ws.Cells.Item(7, 2).Value = -0.04382090528107291
$ws.Cells.Item(7, 3).Value = 1.245761672062962
$ws.Cells.Item(7, 4).Value = -3.257358529865598
$ws.Cells.Item(7, 5).Value = 4.295561381434368
$ws.Cells.Item(7, 6).Value = 0.5488757090113225
$ws.Cells.Item(7, 7).Value = 0.9480759050694929
$ws.Cells.Item(7, 8).Value = 0.00642789908831664
$ws.Cells.Item(7, 9).Value = [double]"6.223607749267224e-06"
$ws.Cells.Item(7, 10).Value = 0.742474602301843
$ws.Cells.Item(7, 11).Value = [double]"1.417294209166032e-27"
$ws.Cells.Item(7, 12).Value = 0.06554819525478191
$ws.Cells.Item(7, 13).Value = 0.2599715093565305
$ws.Cells.Item(7, 14).Value = 0.229143701533319
$ws.Cells.Item(7, 15).Value = 0.2949467309154264
$ws.Cells.Item(8, 2).Value = -0.03330052884810956
$ws.Cells.Item(8, 3).Value = 1.359763759955901
$ws.Cells.Item(8, 4).Value = -5.863022169465559
$ws.Cells.Item(8, 5).Value = 3.10216798942466
$ws.Cells.Item(8, 6).Value = -0.6320568296909684
$ws.Cells.Item(8, 7).Value = 2.303863619888101
$ws.Cells.Item(8, 8).Value = 0.06523960686139221
$ws.Cells.Item(8, 9).Value = [double]"5.624204735266147e-08"
$ws.Cells.Item(8, 10).Value = 1.0884540982985
$ws.Cells.Item(8, 11).Value = [double]"2.932689677292757e-14"
$ws.Cells.Item(8, 12).Value = 0.006982132349433313
$ws.Cells.Item(8, 13).Value = 0.7134317257705455
$ws.Cells.Item(8, 14).Value = 0.6718128544999703
$ws.Cells.Item(8, 15).Value = 0.7576288901390191
$ws.Cells.Item(9, 2).Value = -0.06595281813746115
$ws.Cells.Item(9, 3).Value = 1.831368557925425
$ws.Cells.Item(9, 4).Value = -7.068398169649039
$ws.Cells.Item(9, 5).Value = 4.840078570390868
$ws.Cells.Item(9, 6).Value = -0.855595679574217
$ws.Cells.Item(9, 7).Value = 1.735614844676266
$ws.Cells.Item(9, 8).Value = 0.0001729547387266205
$ws.Cells.Item(9, 9).Value = 0
$ws.Cells.Item(9, 10).Value = 0.9520346051278211
$ws.Cells.Item(9, 11).Value = [double]"1.854166338275856e-15"
$ws.Cells.Item(9, 12).Value = [double]"8.76234337604501e-05"
$ws.Cells.Item(9, 13).Value = 0.5950335802432145
$ws.Cells.Item(9, 14).Value = 0.5565429215192903
$ws.Cells.Item(9, 15).Value = 0.6361862633173132
$ws.Cells.Item(10, 2).Value = -0.02677050590185231
$ws.Cells.Item(10, 3).Value = 1.175744473653722
$ws.Cells.Item(10, 4).Value = -5.601956706765
$ws.Cells.Item(10, 5).Value = 3.218016953682232
$ws.Cells.Item(10, 6).Value = -0.6123506091428452
$ws.Cells.Item(10, 7).Value = 4.010995437566506
$ws.Cells.Item(10, 8).Value = 0.04049500744373022
$ws.Cells.Item(10, 9).Value = 0.001998971789768156
$ws.Cells.Item(10, 10).Value = 1.147256669079903
$ws.Cells.Item(10, 11).Value = [double]"6.579564635074442e-06"
$ws.Cells.Item(10, 12).Value = 0.0001437298303894966
$ws.Cells.Item(10, 13).Value = 0.2629923360687634
$ws.Cells.Item(10, 14).Value = 0.2464451615116112
$ws.Cells.Item(10, 15).Value = 0.2806505447567763
$ws.Cells.Item(11, 2).Value = -0.02654094308338512
$ws.Cells.Item(11, 3).Value = 1.318213528485582
$ws.Cells.Item(11, 4).Value = -4.988443605271979
$ws.Cells.Item(11, 5).Value = 5.156645897593171
$ws.Cells.Item(11, 6).Value = 0.2406932374477504
$ws.Cells.Item(11, 7).Value = 3.206421518147792
$ws.Cells.Item(11, 8).Value = 0.003176824354742507
$ws.Cells.Item(11, 9).Value = [double]"1.429115904505451e-07"
$ws.Cells.Item(11, 10).Value = 1.065811796985765
$ws.Cells.Item(11, 11).Value = 0.0002548176918218574
$ws.Cells.Item(11, 12).Value = 0.0004674627214102582
$ws.Cells.Item(11, 13).Value = 0.06871412692761329
$ws.Cells.Item(11, 14).Value = 0.06245224486320335
$ws.Cells.Item(11, 15).Value = 0.07560386739926653
$ws.Cells.Item(12, 2).Value = -0.01655457148027956
$ws.Cells.Item(12, 3).Value = 1.102208028302038
$ws.Cells.Item(12, 4).Value = -4.465965696818451
$ws.Cells.Item(12, 5).Value = 3.155128727873384
$ws.Cells.Item(12, 6).Value = -0.2592790716929455
$ws.Cells.Item(12, 7).Value = 1.487900519459402
$ws.Cells.Item(12, 8).Value = 0.0200599048354847
$ws.Cells.Item(12, 9).Value = 0.02303905703844644
$ws.Cells.Item(12, 10).Value = 1.249908920067508
$ws.Cells.Item(12, 11).Value = 0.0001857687352268457
$ws.Cells.Item(12, 12).Value = 0.01368010387716044
$ws.Cells.Item(12, 13).Value = 0.4636733120359737
$ws.Cells.Item(12, 14).Value = 0.4457127095215158
$ws.Cells.Item(12, 15).Value = 0.4823576615645759
$ws.Cells.Item(13, 2).Value = -0.02148595520095248
$ws.Cells.Item(13, 3).Value = 1.106062008435703
$ws.Cells.Item(13, 4).Value = -4.805809207206539
$ws.Cells.Item(13, 5).Value = 2.170405330564763
$ws.Cells.Item(13, 6).Value = -0.4553036466822096
$ws.Cells.Item(13, 7).Value = 1.714041975761713
$ws.Cells.Item(13, 8).Value = 0.03953950434075763
$ws.Cells.Item(13, 9).Value = 0.02390773039945093
$ws.Cells.Item(13, 10).Value = 0.948135836005815
$ws.Cells.Item(13, 11).Value = [double]"7.389524213593849e-07"
$ws.Cells.Item(13, 12).Value = 0.000823339131625459
$ws.Cells.Item(13, 13).Value = 3.664655073384121
$ws.Cells.Item(13, 14).Value = 3.503387727080372
$ws.Cells.Item(13, 15).Value = 3.83334585066664
$ws.Cells.Item(14, 2).Value = -0.01880754928364788
$ws.Cells.Item(14, 3).Value = 1.095059319020325
$ws.Cells.Item(14, 4).Value = -4.475344338412787
$ws.Cells.Item(14, 5).Value = 3.875840694318431
$ws.Cells.Item(14, 6).Value = -0.03663271877713174
$ws.Cells.Item(14, 7).Value = 2.622145387879804
$ws.Cells.Item(14, 8).Value = 0.01243312250406129
$ws.Cells.Item(14, 9).Value = 0.02528631327847874
$ws.Cells.Item(14, 10).Value = 1.237297169901697
$ws.Cells.Item(14, 11).Value = [double]"2.310954593049932e-08"
$ws.Cells.Item(14, 12).Value = 0.006595937298413504
$ws.Cells.Item(14, 13).Value = 4.301765034404673
$ws.Cells.Item(14, 14).Value = 4.184095796292002
$ws.Cells.Item(14, 15).Value = 4.422743482026907
$ws.Cells.Item(15, 2).Value = -0.03105870820959822
$ws.Cells.Item(15, 3).Value = 1.301421231963948
$ws.Cells.Item(15, 4).Value = -2.75757037817916
$ws.Cells.Item(15, 5).Value = 5.999898454250217
$ws.Cells.Item(15, 6).Value = 1.032388797505195
$ws.Cells.Item(15, 7).Value = 3.143364563972882
$ws.Cells.Item(15, 8).Value = 0.1070792365366116
$ws.Cells.Item(15, 9).Value = [double]"8.113772687057264e-08"
$ws.Cells.Item(15, 10).Value = 0.884457365323939
$ws.Cells.Item(15, 11).Value = 0.0005620887572302592
$ws.Cells.Item(15, 12).Value = 0.0004167539380658026
$ws.Cells.Item(15, 13).Value = 0.8662714033921668
$ws.Cells.Item(15, 14).Value = 0.8027575132605852
$ws.Cells.Item(15, 15).Value = 0.9348104900158516
$ws.Cells.Item(16, 2).Value = -0.0171220513330381
$ws.Cells.Item(16, 3).Value = 1.032575113670937
$ws.Cells.Item(16, 4).Value = -4.931772645592038
$ws.Cells.Item(16, 5).Value = 2.438457937112157
$ws.Cells.Item(16, 6).Value = -0.6846942161887539
$ws.Cells.Item(16, 7).Value = 3.282576437445911
$ws.Cells.Item(16, 8).Value = 0.005067020828246594
$ws.Cells.Item(16, 9).Value = 0.1825950782664085
$ws.Cells.Item(16, 10).Value = 1.302130080581972
$ws.Cells.Item(16, 11).Value = 0.0003640059742283136
$ws.Cells.Item(16, 12).Value = 0.0002378665670032319
$ws.Cells.Item(16, 13).Value = 0.606043049090066
$ws.Cells.Item(16, 14).Value = 0.5888667861975406
$ws.Cells.Item(16, 15).Value = 0.6237203149494221
$ws.Cells.Item(17, 2).Value = -0.01826864653641791
$ws.Cells.Item(17, 3).Value = 1.090383743179026
$ws.Cells.Item(17, 4).Value = -5.088350891306335
$ws.Cells.Item(17, 5).Value = 2.937380873793336
$ws.Cells.Item(17, 6).Value = -0.4898662704049352
$ws.Cells.Item(17, 7).Value = 3.262685697802854
$ws.Cells.Item(17, 8).Value = 0.02296910544022162
$ws.Cells.Item(17, 9).Value = 0.0381183193857223
$ws.Cells.Item(17, 10).Value = 1.14421173473669
$ws.Cells.Item(17, 11).Value = [double]"6.354233422379444e-06"
$ws.Cells.Item(17, 12).Value = 0.0003947229899670635
$ws.Cells.Item(17, 13).Value = 11.73938767699414
$ws.Cells.Item(17, 14).Value = 11.37848929732894
$ws.Cells.Item(17, 15).Value = 12.111732887345
